{"js": "/*\n * Replaces each arithmetic-expression cell's text with its updated value.\n * The document is a 5-column table of \"NN+NN=\" / \"NN-NN=\" style problems;\n * this script performs an exact (non-wildcard) search for each original\n * expression and replaces the first match in place via insertText, which\n * preserves the existing run formatting (font/size) of the matched range.\n */\nconst replacements = [{\"o\": \"21+9=\", \"n\": \"73-73=\"}, {\"o\": \"3+87=\", \"n\": \"51+27=\"}, {\"o\": \"57-30=\", \"n\": \"51-22=\"}, {\"o\": \"54-16=\", \"n\": \"86-74=\"}, {\"o\": \"71+20=\", \"n\": \"35+34=\"}, {\"o\": \"48-22=\", \"n\": \"80-5=\"}, {\"o\": \"26+17=\", \"n\": \"31+37=\"}, {\"o\": \"39-21=\", \"n\": \"90-32=\"}, {\"o\": \"76-56=\", \"n\": \"69+12=\"}, {\"o\": \"22+24=\", \"n\": \"64+26=\"}, {\"o\": \"84-6=\", \"n\": \"28+30=\"}, {\"o\": \"80-64=\", \"n\": \"28+15=\"}, {\"o\": \"98-48=\", \"n\": \"83-19=\"}, {\"o\": \"71+13=\", \"n\": \"71-13=\"}, {\"o\": \"60-20=\", \"n\": \"21+75=\"}, {\"o\": \"82-28=\", \"n\": \"77-61=\"}, {\"o\": \"18-15=\", \"n\": \"76+12=\"}, {\"o\": \"55-22=\", \"n\": \"27+37=\"}, {\"o\": \"88+5=\", \"n\": \"10+84=\"}, {\"o\": \"0+96=\", \"n\": \"83-69=\"}, {\"o\": \"78-75=\", \"n\": \"87-70=\"}, {\"o\": \"93-43=\", \"n\": \"91-60=\"}, {\"o\": \"4+90=\", \"n\": \"84-69=\"}, {\"o\": \"97-29=\", \"n\": \"76-53=\"}, {\"o\": \"5+39=\", \"n\": \"97-78=\"}, {\"o\": \"66+3=\", \"n\": \"2+44=\"}, {\"o\": \"90-15=\", \"n\": \"12+52=\"}, {\"o\": \"17-3=\", \"n\": \"78-54=\"}, {\"o\": \"26+48=\", \"n\": \"20+0=\"}, {\"o\": \"27-1=\", \"n\": \"22+15=\"}, {\"o\": \"14-12=\", \"n\": \"57-41=\"}, {\"o\": \"71+17=\", \"n\": \"23+25=\"}, {\"o\": \"67-10=\", \"n\": \"72-68=\"}, {\"o\": \"98-36=\", \"n\": \"5+3=\"}, {\"o\": \"39-25=\", \"n\": \"22+12=\"}, {\"o\": \"86-7=\", \"n\": \"74-2=\"}, {\"o\": \"86-85=\", \"n\": \"48-42=\"}, {\"o\": \"21-1=\", \"n\": \"67-58=\"}, {\"o\": \"46-20=\", \"n\": \"95-49=\"}, {\"o\": \"43-2=\", \"n\": \"92-18=\"}, {\"o\": \"93-4=\", \"n\": \"30-24=\"}, {\"o\": \"53+23=\", \"n\": \"70+27=\"}, {\"o\": \"88-50=\", \"n\": \"22+5=\"}, {\"o\": \"42+36=\", \"n\": \"89-74=\"}, {\"o\": \"51+36=\", \"n\": \"69-60=\"}, {\"o\": \"73-4=\", \"n\": \"28-19=\"}, {\"o\": \"0+20=\", \"n\": \"39+10=\"}, {\"o\": \"51-6=\", \"n\": \"54-5=\"}, {\"o\": \"48-24=\", \"n\": \"9+11=\"}, {\"o\": \"41+35=\", \"n\": \"78-39=\"}, {\"o\": \"67-20=\", \"n\": \"10+76=\"}, {\"o\": \"52+35=\", \"n\": \"19+77=\"}, {\"o\": \"92-16=\", \"n\": \"19+10=\"}, {\"o\": \"38+2=\", \"n\": \"1+91=\"}, {\"o\": \"30-2=\", \"n\": \"62+33=\"}, {\"o\": \"72-3=\", \"n\": \"90-51=\"}, {\"o\": \"80-3=\", \"n\": \"70-47=\"}, {\"o\": \"72-51=\", \"n\": \"16+59=\"}, {\"o\": \"5+62=\", \"n\": \"65-46=\"}, {\"o\": \"51-38=\", \"n\": \"61-3=\"}, {\"o\": \"0+28=\", \"n\": \"79-73=\"}, {\"o\": \"32-18=\", \"n\": \"65+4=\"}, {\"o\": \"8+88=\", \"n\": \"98-92=\"}, {\"o\": \"9+68=\", \"n\": \"52+8=\"}, {\"o\": \"41-11=\", \"n\": \"76-11=\"}, {\"o\": \"82-3=\", \"n\": \"93-67=\"}, {\"o\": \"20+38=\", \"n\": \"18+45=\"}, {\"o\": \"77-11=\", \"n\": \"50+2=\"}, {\"o\": \"30+59=\", \"n\": \"12+82=\"}, {\"o\": \"29+28=\", \"n\": \"88-30=\"}, {\"o\": \"34-17=\", \"n\": \"70-40=\"}, {\"o\": \"41+43=\", \"n\": \"71-66=\"}, {\"o\": \"15+63=\", \"n\": \"78-25=\"}, {\"o\": \"63+14=\", \"n\": \"26+37=\"}, {\"o\": \"30+10=\", \"n\": \"76-26=\"}, {\"o\": \"60-15=\", \"n\": \"3+89=\"}, {\"o\": \"23+32=\", \"n\": \"3+7=\"}, {\"o\": \"31+6=\", \"n\": \"75-26=\"}, {\"o\": \"37+27=\", \"n\": \"84+2=\"}, {\"o\": \"5+58=\", \"n\": \"37+8=\"}, {\"o\": \"85-76=\", \"n\": \"73-69=\"}, {\"o\": \"82-76=\", \"n\": \"61-22=\"}, {\"o\": \"13+6=\", \"n\": \"48+48=\"}, {\"o\": \"78-58=\", \"n\": \"48-47=\"}, {\"o\": \"24+57=\", \"n\": \"64-43=\"}, {\"o\": \"58-6=\", \"n\": \"84-41=\"}, {\"o\": \"19+66=\", \"n\": \"99-7=\"}, {\"o\": \"66+2=\", \"n\": \"49-1=\"}, {\"o\": \"16+83=\", \"n\": \"84+0=\"}, {\"o\": \"26+58=\", \"n\": \"30+15=\"}, {\"o\": \"38+43=\", \"n\": \"50+10=\"}, {\"o\": \"28-8=\", \"n\": \"22+69=\"}, {\"o\": \"16-2=\", \"n\": \"61+11=\"}, {\"o\": \"63-36=\", \"n\": \"80-80=\"}, {\"o\": \"88-12=\", \"n\": \"55-8=\"}, {\"o\": \"80-46=\", \"n\": \"3+26=\"}, {\"o\": \"18+52=\", \"n\": \"27+34=\"}, {\"o\": \"23+2=\", \"n\": \"49-10=\"}, {\"o\": \"29+25=\", \"n\": \"35+57=\"}, {\"o\": \"18+16=\", \"n\": \"35-6=\"}];\n\nfor (const {o, n} of replacements) {\n  const results = context.document.body.search(o, {matchCase: true, matchWildcards: false});\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + o);\n  }\n\n  results.items[0].insertText(n, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Updates each arithmetic-expression cell in the table with its new value.\n# Uses Find/Execute with wdReplaceAll semantics on a full-document range for\n# each unique original expression; Find.Execute preserves the formatting of\n# the matched text because it replaces only the matched characters in place.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"21+9=\", \"73-73=\"),\n    @(\"3+87=\", \"51+27=\"),\n    @(\"57-30=\", \"51-22=\"),\n    @(\"54-16=\", \"86-74=\"),\n    @(\"71+20=\", \"35+34=\"),\n    @(\"48-22=\", \"80-5=\"),\n    @(\"26+17=\", \"31+37=\"),\n    @(\"39-21=\", \"90-32=\"),\n    @(\"76-56=\", \"69+12=\"),\n    @(\"22+24=\", \"64+26=\"),\n    @(\"84-6=\", \"28+30=\"),\n    @(\"80-64=\", \"28+15=\"),\n    @(\"98-48=\", \"83-19=\"),\n    @(\"71+13=\", \"71-13=\"),\n    @(\"60-20=\", \"21+75=\"),\n    @(\"82-28=\", \"77-61=\"),\n    @(\"18-15=\", \"76+12=\"),\n    @(\"55-22=\", \"27+37=\"),\n    @(\"88+5=\", \"10+84=\"),\n    @(\"0+96=\", \"83-69=\"),\n    @(\"78-75=\", \"87-70=\"),\n    @(\"93-43=\", \"91-60=\"),\n    @(\"4+90=\", \"84-69=\"),\n    @(\"97-29=\", \"76-53=\"),\n    @(\"5+39=\", \"97-78=\"),\n    @(\"66+3=\", \"2+44=\"),\n    @(\"90-15=\", \"12+52=\"),\n    @(\"17-3=\", \"78-54=\"),\n    @(\"26+48=\", \"20+0=\"),\n    @(\"27-1=\", \"22+15=\"),\n    @(\"14-12=\", \"57-41=\"),\n    @(\"71+17=\", \"23+25=\"),\n    @(\"67-10=\", \"72-68=\"),\n    @(\"98-36=\", \"5+3=\"),\n    @(\"39-25=\", \"22+12=\"),\n    @(\"86-7=\", \"74-2=\"),\n    @(\"86-85=\", \"48-42=\"),\n    @(\"21-1=\", \"67-58=\"),\n    @(\"46-20=\", \"95-49=\"),\n    @(\"43-2=\", \"92-18=\"),\n    @(\"93-4=\", \"30-24=\"),\n    @(\"53+23=\", \"70+27=\"),\n    @(\"88-50=\", \"22+5=\"),\n    @(\"42+36=\", \"89-74=\"),\n    @(\"51+36=\", \"69-60=\"),\n    @(\"73-4=\", \"28-19=\"),\n    @(\"0+20=\", \"39+10=\"),\n    @(\"51-6=\", \"54-5=\"),\n    @(\"48-24=\", \"9+11=\"),\n    @(\"41+35=\", \"78-39=\"),\n    @(\"67-20=\", \"10+76=\"),\n    @(\"52+35=\", \"19+77=\"),\n    @(\"92-16=\", \"19+10=\"),\n    @(\"38+2=\", \"1+91=\"),\n    @(\"30-2=\", \"62+33=\"),\n    @(\"72-3=\", \"90-51=\"),\n    @(\"80-3=\", \"70-47=\"),\n    @(\"72-51=\", \"16+59=\"),\n    @(\"5+62=\", \"65-46=\"),\n    @(\"51-38=\", \"61-3=\"),\n    @(\"0+28=\", \"79-73=\"),\n    @(\"32-18=\", \"65+4=\"),\n    @(\"8+88=\", \"98-92=\"),\n    @(\"9+68=\", \"52+8=\"),\n    @(\"41-11=\", \"76-11=\"),\n    @(\"82-3=\", \"93-67=\"),\n    @(\"20+38=\", \"18+45=\"),\n    @(\"77-11=\", \"50+2=\"),\n    @(\"30+59=\", \"12+82=\"),\n    @(\"29+28=\", \"88-30=\"),\n    @(\"34-17=\", \"70-40=\"),\n    @(\"41+43=\", \"71-66=\"),\n    @(\"15+63=\", \"78-25=\"),\n    @(\"63+14=\", \"26+37=\"),\n    @(\"30+10=\", \"76-26=\"),\n    @(\"60-15=\", \"3+89=\"),\n    @(\"23+32=\", \"3+7=\"),\n    @(\"31+6=\", \"75-26=\"),\n    @(\"37+27=\", \"84+2=\"),\n    @(\"5+58=\", \"37+8=\"),\n    @(\"85-76=\", \"73-69=\"),\n    @(\"82-76=\", \"61-22=\"),\n    @(\"13+6=\", \"48+48=\"),\n    @(\"78-58=\", \"48-47=\"),\n    @(\"24+57=\", \"64-43=\"),\n    @(\"58-6=\", \"84-41=\"),\n    @(\"19+66=\", \"99-7=\"),\n    @(\"66+2=\", \"49-1=\"),\n    @(\"16+83=\", \"84+0=\"),\n    @(\"26+58=\", \"30+15=\"),\n    @(\"38+43=\", \"50+10=\"),\n    @(\"28-8=\", \"22+69=\"),\n    @(\"16-2=\", \"61+11=\"),\n    @(\"63-36=\", \"80-80=\"),\n    @(\"88-12=\", \"55-8=\"),\n    @(\"80-46=\", \"3+26=\"),\n    @(\"18+52=\", \"27+34=\"),\n    @(\"23+2=\", \"49-10=\"),\n    @(\"29+25=\", \"35+57=\"),\n    @(\"18+16=\", \"35-6=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"No match found for: $old\"\n    }\n}\n"}
